# Refine metadata to be an additional tab:
#  1. Add a new "metadata" worksheet after "data" with summary info about
#     the panel query (name/id/version/timestamps/request URL).
#  2. Refresh the "time_taken" (column F) timestamps on the "data" sheet to
#     reflect the new query run.

$wb = $excel.ActiveWorkbook
$dataWs = $wb.Worksheets.Item("data")

# ---------------------------------------------------------------------------
# 1. Update the per-row query timestamps on the "data" sheet (column F).
# ---------------------------------------------------------------------------
$timeTaken = @{
    "2" = "2021-10-05 14:33:27.517665";
    "3" = "2021-10-05 14:33:27.517671";
    "4" = "2021-10-05 14:33:27.517674";
    "5" = "2021-10-05 14:33:27.517676";
    "6" = "2021-10-05 14:33:27.517678";
    "7" = "2021-10-05 14:33:27.517680";
    "8" = "2021-10-05 14:33:27.517682";
    "9" = "2021-10-05 14:33:27.517684";
    "10" = "2021-10-05 14:33:27.517686";
    "11" = "2021-10-05 14:33:27.517688";
    "12" = "2021-10-05 14:33:27.517690";
    "13" = "2021-10-05 14:33:27.517691";
    "14" = "2021-10-05 14:33:27.517693";
    "15" = "2021-10-05 14:33:27.517695";
    "16" = "2021-10-05 14:33:27.517697";
    "17" = "2021-10-05 14:33:27.517699";
    "18" = "2021-10-05 14:33:27.517701";
    "19" = "2021-10-05 14:33:27.517703";
    "20" = "2021-10-05 14:33:27.517705";
    "21" = "2021-10-05 14:33:27.517707";
    "22" = "2021-10-05 14:33:27.517709";
    "23" = "2021-10-05 14:33:27.517711";
    "24" = "2021-10-05 14:33:27.517713";
    "25" = "2021-10-05 14:33:27.517715";
    "26" = "2021-10-05 14:33:27.517717";
    "27" = "2021-10-05 14:33:27.517719";
    "28" = "2021-10-05 14:33:27.517721";
    "29" = "2021-10-05 14:33:27.517723";
    "30" = "2021-10-05 14:33:27.517724";
    "31" = "2021-10-05 14:33:27.517726";
    "32" = "2021-10-05 14:33:27.517728";
    "33" = "2021-10-05 14:33:27.517730";
    "34" = "2021-10-05 14:33:27.517732";
    "35" = "2021-10-05 14:33:27.517734";
    "36" = "2021-10-05 14:33:27.517736";
    "37" = "2021-10-05 14:33:27.517738";
    "38" = "2021-10-05 14:33:27.517740";
    "39" = "2021-10-05 14:33:27.517741";
    "40" = "2021-10-05 14:33:27.517743";
    "41" = "2021-10-05 14:33:27.517745";
    "42" = "2021-10-05 14:33:27.517747";
    "43" = "2021-10-05 14:33:27.517749";
    "44" = "2021-10-05 14:33:27.517751";
    "45" = "2021-10-05 14:33:27.517753";
    "46" = "2021-10-05 14:33:27.517755";
    "47" = "2021-10-05 14:33:27.517756";
    "48" = "2021-10-05 14:33:27.517758";
    "49" = "2021-10-05 14:33:27.517760";
    "50" = "2021-10-05 14:33:27.517762";
    "51" = "2021-10-05 14:33:27.517764";
    "52" = "2021-10-05 14:33:27.517766";
    "53" = "2021-10-05 14:33:27.517768";
    "54" = "2021-10-05 14:33:27.517770";
    "55" = "2021-10-05 14:33:27.517772";
    "56" = "2021-10-05 14:33:27.517774";
    "57" = "2021-10-05 14:33:27.517775"
}

foreach ($row in $timeTaken.Keys) {
    $dataWs.Range("F" + $row).Value = $timeTaken[$row]
}

# ---------------------------------------------------------------------------
# 2. Add the new "metadata" worksheet right after "data".
# ---------------------------------------------------------------------------
$metaWs = $wb.Worksheets.Add($null, $dataWs)
$metaWs.Name = "metadata"

# Match the header / A-column formatting used on the "data" sheet (bold,
# bordered, centered) by copying the existing formatted cells rather than
# re-building the style from scratch.
$dataWs.Range("A2").Copy()
$metaWs.Range("A2").PasteSpecial(-4122)

$dataWs.Range("B1:F1").Copy()
$metaWs.Range("B1:G1").PasteSpecial(-4122)

# Header row.
$metaWs.Range("B1").Value = "data_name"
$metaWs.Range("C1").Value = "data_id"
$metaWs.Range("D1").Value = "data_version"
$metaWs.Range("E1").Value = "data_version_created"
$metaWs.Range("F1").Value = "panel_query_time"
$metaWs.Range("G1").Value = "panel_get_request"

# Data row.
$metaWs.Range("A2").Value = 0
$metaWs.Range("B2").Value = "Chromosome Breakage Disorders"
$metaWs.Range("C2").Value = 79
$metaWs.Range("D2").NumberFormat = "@"
$metaWs.Range("D2").Value = "1.5"
$metaWs.Range("E2").Value = "2021-10-02T22:58:40.376976Z"
$metaWs.Range("F2").Value = "2021-10-05 14:33:27.515125"
$metaWs.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/79/?format=json"

# Keep "data" as the active/selected sheet (unchanged bookViews in the diff).
$dataWs.Activate()
[void]$dataWs.Range("A1").Select()
